# chore: update Sheets via scheduled runner
# Refreshes currentAveragePrice/NQ/HQ, LevePrice and LeveProfit columns
# (H:N) for a number of leve rows across the ALC, ARM, BSM, CRP, CUL,
# GSM and WVR sheets, reflecting the latest market-board figures.

$wb = $excel.ActiveWorkbook

# Hunk 0: ALC!row53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 215.2
$ws.Range("I53").Value = 246.25
$ws.Range("J53").Value = 194.5
$ws.Range("K53").Value = 246.25
$ws.Range("L53").Value = 194.5
$ws.Range("M53").Value = 390.75
$ws.Range("N53").Value = -1468.5

# Hunk 1: ALC!row64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3723.3948
$ws.Range("I64").Value = 3594.762
$ws.Range("J64").Value = 3882.2942
$ws.Range("K64").Value = 3594.762
$ws.Range("L64").Value = 3882.2942
$ws.Range("M64").Value = -3346.762
$ws.Range("N64").Value = -4378.2942

# Hunk 2: ALC!row67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3723.3948
$ws.Range("I67").Value = 3594.762
$ws.Range("J67").Value = 3882.2942
$ws.Range("K67").Value = 3594.762
$ws.Range("L67").Value = 3882.2942
$ws.Range("M67").Value = -2736.762
$ws.Range("N67").Value = -5598.2942

# Hunk 3: ALC!row99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 375.25
$ws.Range("I99").Value = 302.5
$ws.Range("J99").Value = 496.5
$ws.Range("K99").Value = 907.5
$ws.Range("L99").Value = 1489.5
$ws.Range("M99").Value = 590.5
$ws.Range("N99").Value = -4485.5

# Hunk 4: ALC!row125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2600
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2600
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 23400
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -28320

# Hunk 5: ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 22036.92
$ws.Range("I137").Value = 1266.2122
$ws.Range("J137").Value = 62356.53
$ws.Range("K137").Value = 3798.6366
$ws.Range("L137").Value = 187069.59
$ws.Range("M137").Value = -1248.6366
$ws.Range("N137").Value = -192169.59

# Hunk 6: ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3273.7258
$ws.Range("I138").Value = 1182.3334
$ws.Range("J138").Value = 4129.2954
$ws.Range("K138").Value = 3547.0002
$ws.Range("L138").Value = 12387.8862
$ws.Range("M138").Value = 1592.9998
$ws.Range("N138").Value = -22667.8862

# Hunk 7: ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16671128
$ws.Range("I32").Value = 18521476
$ws.Range("J32").Value = 18002.166
$ws.Range("K32").Value = 18521476
$ws.Range("L32").Value = 18002.166
$ws.Range("M32").Value = -18521189
$ws.Range("N32").Value = -18576.166

# Hunk 8: ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1838.7
$ws.Range("I45").Value = 1029.8334
$ws.Range("J45").Value = 3052
$ws.Range("K45").Value = 1029.8334
$ws.Range("L45").Value = 3052
$ws.Range("M45").Value = -652.8334
$ws.Range("N45").Value = -3806

# Hunk 9: BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

# Hunk 10: BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

# Hunk 11: BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 509.57144
$ws.Range("I94").Value = 530.9
$ws.Range("J94").Value = 456.25
$ws.Range("K94").Value = 530.9
$ws.Range("L94").Value = 456.25
$ws.Range("M94").Value = -79.89999999999998
$ws.Range("N94").Value = -1358.25

# Hunk 12: BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2944.94
$ws.Range("I105").Value = 1817.5
$ws.Range("J105").Value = 2991.9167
$ws.Range("K105").Value = 1817.5
$ws.Range("L105").Value = 2991.9167
$ws.Range("M105").Value = -70.5
$ws.Range("N105").Value = -6485.9167

# Hunk 13: BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1903.08
$ws.Range("I134").Value = 1894.6522
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 5683.9566
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -3148.9566
$ws.Range("N134").Value = -11070

# Hunk 14: CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 753.92755
$ws.Range("I58").Value = 670.2982
$ws.Range("J58").Value = 1151.1666
$ws.Range("K58").Value = 670.2982
$ws.Range("L58").Value = 1151.1666
$ws.Range("M58").Value = -467.2982
$ws.Range("N58").Value = -1557.1666

# Hunk 15: CRP!row86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4539.125
$ws.Range("I86").Value = 6258.143
$ws.Range("J86").Value = 3202.111
$ws.Range("K86").Value = 6258.143
$ws.Range("L86").Value = 3202.111
$ws.Range("M86").Value = -5135.143
$ws.Range("N86").Value = -5448.111

# Hunk 16: CRP!row89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 4539.125
$ws.Range("I89").Value = 6258.143
$ws.Range("J89").Value = 3202.111
$ws.Range("K89").Value = 31290.715
$ws.Range("L89").Value = 16010.555
$ws.Range("M89").Value = -25674.715
$ws.Range("N89").Value = -27242.555

# Hunk 17: CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 753.92755
$ws.Range("I136").Value = 670.2982
$ws.Range("J136").Value = 1151.1666
$ws.Range("K136").Value = 2010.8946
$ws.Range("L136").Value = 3453.4998
$ws.Range("M136").Value = 539.1054000000001
$ws.Range("N136").Value = -8553.4998

# Hunk 18: CUL!row86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1032.409
$ws.Range("I86").Value = 803.3333
$ws.Range("J86").Value = 1191
$ws.Range("K86").Value = 2409.9999
$ws.Range("L86").Value = 3573
$ws.Range("M86").Value = -1223.9999
$ws.Range("N86").Value = -5945

# Hunk 19: CUL!row89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1032.409
$ws.Range("I89").Value = 803.3333
$ws.Range("J89").Value = 1191
$ws.Range("K89").Value = 7229.9997
$ws.Range("L89").Value = 10719
$ws.Range("M89").Value = -1301.9997
$ws.Range("N89").Value = -22575

# Hunk 20: CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 744.8868
$ws.Range("I131").Value = 453.22223
$ws.Range("J131").Value = 894.8857400000001
$ws.Range("K131").Value = 1359.66669
$ws.Range("L131").Value = 2684.65722
$ws.Range("M131").Value = 3680.33331
$ws.Range("N131").Value = -12764.65722

# Hunk 21: CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 722885.5600000001
$ws.Range("I132").Value = 1222.2222
$ws.Range("J132").Value = 2021879.6
$ws.Range("K132").Value = 10999.9998
$ws.Range("L132").Value = 18196916.4
$ws.Range("M132").Value = -8469.9998
$ws.Range("N132").Value = -18201976.4

# Hunk 22: GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7167.696
$ws.Range("I70").Value = 7826.353
$ws.Range("J70").Value = 5301.5
$ws.Range("K70").Value = 7826.353
$ws.Range("L70").Value = 5301.5
$ws.Range("M70").Value = -7556.353
$ws.Range("N70").Value = -5841.5

# Hunk 23: GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7167.696
$ws.Range("I73").Value = 7826.353
$ws.Range("J73").Value = 5301.5
$ws.Range("K73").Value = 7826.353
$ws.Range("L73").Value = 5301.5
$ws.Range("M73").Value = -6890.353
$ws.Range("N73").Value = -7173.5

# Hunk 24: GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1896.2307
$ws.Range("I122").Value = 1629
$ws.Range("J122").Value = 2497.5
$ws.Range("K122").Value = 4887
$ws.Range("L122").Value = 7492.5
$ws.Range("M122").Value = -2437
$ws.Range("N122").Value = -12392.5

# Hunk 25: GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4693.7
$ws.Range("I132").Value = 4740.5
$ws.Range("J132").Value = 4506.5
$ws.Range("K132").Value = 14221.5
$ws.Range("L132").Value = 13519.5
$ws.Range("M132").Value = -11691.5
$ws.Range("N132").Value = -18579.5

# Hunk 26: WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3475.5
$ws.Range("I126").Value = 4333.3335
$ws.Range("J126").Value = 3107.8572
$ws.Range("K126").Value = 13000.0005
$ws.Range("L126").Value = 9323.571599999999
$ws.Range("M126").Value = -10530.0005
$ws.Range("N126").Value = -14263.5716

# Hunk 27: WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1025.1167
$ws.Range("I132").Value = 674.55316
$ws.Range("J132").Value = 2292.5386
$ws.Range("K132").Value = 2023.65948
$ws.Range("L132").Value = 6877.6158
$ws.Range("M132").Value = 506.3405199999997
$ws.Range("N132").Value = -11937.6158

# Hunk 28: WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5538.478
$ws.Range("I136").Value = 6249.1665
$ws.Range("J136").Value = 2980
$ws.Range("K136").Value = 18747.4995
$ws.Range("L136").Value = 8940
$ws.Range("M136").Value = -16197.4995
$ws.Range("N136").Value = -14040
